$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must remain stored as text
# (matching the workbook's existing convention of inlineStr/text price cells).
# Force text format on the affected D cells before assigning, otherwise the
# engine will auto-convert numeric-looking strings into real numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values, row by row, matching the published diff.
$ws.Range("D2").Value = "244.49"
$ws.Range("D3").Value = "23.69"
$ws.Range("D4").Value = "5.193"
$ws.Range("D5").Value = "0.05752"
$ws.Range("D6").Value = "6.450"
$ws.Range("D7").Value = "3.255"
$ws.Range("D8").Value = "0.8096"
$ws.Range("D9").Value = "0.8668"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01010"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1375"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.06932"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03172"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03046"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09328"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.813"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001522"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04712"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "0.006157"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "0.001235"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.004088"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "0.00008694"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "3.548"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "2.146"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D26").Value = "0.1295"
$ws.Range("D27").Value = "0.0002327"
$ws.Range("D40").Value = "0.03702"
$ws.Range("D41").Value = "0.006235"
$ws.Range("D42").Value = "0.1049"
$ws.Range("D43").Value = "0.002598"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "0.007084"
$ws.Range("D45").Value = "0.00005326"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").Value = "0.4297"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "0.002048"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("D50").Value = "0.0001999"
